$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 219, shifting rows 219:264 down to 220:265
$ws.Rows.Item(219).Insert()

# Populate the newly inserted row 219 with the new record's data
$ws.Cells.Item(219, 1).Value = 10
$ws.Cells.Item(219, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(219, 3).Value = "La Araucanía"
$ws.Cells.Item(219, 4).Value = 44782
$ws.Cells.Item(219, 5).Value = 9
$ws.Cells.Item(219, 6).Value = 100112043
$ws.Cells.Item(219, 7).Value = "Pepino dulce"
$ws.Cells.Item(219, 8).Value = "Cultivar IV Región"
$ws.Cells.Item(219, 9).Value = "Primera"
$ws.Cells.Item(219, 10).Value = 130
$ws.Cells.Item(219, 11).Value = 18000
$ws.Cells.Item(219, 12).Value = 19000
$ws.Cells.Item(219, 13).Value = 18615
$ws.Cells.Item(219, 14).Value = "$/bandeja 18 kilos"
$ws.Cells.Item(219, 15).Value = "Provincia de Limarí"
$ws.Cells.Item(219, 16).Value = 1034
$ws.Cells.Item(219, 17).Value = 18
$ws.Cells.Item(219, 18).Value = "Hortaliza"
